# Updates the cryptos list with new price/volume figures.
# Generated from the authoritative diff: updates D (Price) and E (Volume(1h))
# columns for the affected rows, writing values as plain text so that
# numeric-looking strings (e.g. trailing zeros, multi-dot prices) survive
# the round-trip exactly as in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '29.337.06' }
    @{ Cell = 'E2'; Value = '  +0.01%  ' }
    @{ Cell = 'D3'; Value = '1.876.78' }
    @{ Cell = 'E3'; Value = '  +0.13%  ' }
    @{ Cell = 'D4'; Value = '1.002' }
    @{ Cell = 'E4'; Value = '  +0.08%  ' }
    @{ Cell = 'D5'; Value = '0.7118' }
    @{ Cell = 'E5'; Value = '  -0.04%  ' }
    @{ Cell = 'D6'; Value = '243.07' }
    @{ Cell = 'E6'; Value = '  +0.47%  ' }
    @{ Cell = 'E8'; Value = '  +2.59%  ' }
    @{ Cell = 'D9'; Value = '0.3155' }
    @{ Cell = 'E9'; Value = '  +1.40%  ' }
    @{ Cell = 'D10'; Value = '25.00' }
    @{ Cell = 'E10'; Value = '  -0.65%  ' }
    @{ Cell = 'D11'; Value = '0.08257' }
    @{ Cell = 'E11'; Value = '  -1.82%  ' }
    @{ Cell = 'D12'; Value = '1.896.36' }
    @{ Cell = 'E12'; Value = '  +0.94%  ' }
    @{ Cell = 'D13'; Value = '5.250' }
    @{ Cell = 'E13'; Value = '  +0.10%  ' }
    @{ Cell = 'D14'; Value = '94.66' }
    @{ Cell = 'E14'; Value = '  +3.81%  ' }
    @{ Cell = 'D15'; Value = '0.7130' }
    @{ Cell = 'E15'; Value = '  +0.15%  ' }
    @{ Cell = 'D16'; Value = '6.377' }
    @{ Cell = 'E16'; Value = '  +4.65%  ' }
    @{ Cell = 'D17'; Value = '0.000008550' }
    @{ Cell = 'E17'; Value = '  +4.05%  ' }
    @{ Cell = 'D18'; Value = '29.348.88' }
    @{ Cell = 'E18'; Value = '  +0.03%  ' }
    @{ Cell = 'D19'; Value = '245.03' }
    @{ Cell = 'D20'; Value = '2.155.14' }
    @{ Cell = 'E20'; Value = '  +1.46%  ' }
    @{ Cell = 'D21'; Value = '13.27' }
    @{ Cell = 'E21'; Value = '  +0.47%  ' }
    @{ Cell = 'E22'; Value = '  +0.03%  ' }
    @{ Cell = 'D23'; Value = '7.787' }
    @{ Cell = 'E23'; Value = '  +0.27%  ' }
    @{ Cell = 'D24'; Value = '1.002' }
    @{ Cell = 'E24'; Value = '  +0.09%  ' }
    @{ Cell = 'D25'; Value = '0.1558' }
    @{ Cell = 'E25'; Value = '  -2.47%  ' }
    @{ Cell = 'D26'; Value = '9.055' }
    @{ Cell = 'E26'; Value = '  +0.21%  ' }
    @{ Cell = 'D27'; Value = '162.61' }
    @{ Cell = 'E27'; Value = '  -0.05%  ' }
    @{ Cell = 'D28'; Value = '18.55' }
    @{ Cell = 'E28'; Value = '  +0.15%  ' }
    @{ Cell = 'D29'; Value = '1.504' }
    @{ Cell = 'E29'; Value = '  -0.41%  ' }
    @{ Cell = 'D30'; Value = '4.422' }
    @{ Cell = 'E30'; Value = '  +0.00%  ' }
    @{ Cell = 'D31'; Value = '4.318' }
    @{ Cell = 'E31'; Value = '  +0.16%  ' }
    @{ Cell = 'D32'; Value = '1.191' }
    @{ Cell = 'E32'; Value = '  -7.75%  ' }
    @{ Cell = 'D33'; Value = '0.05389' }
    @{ Cell = 'E33'; Value = '  +1.83%  ' }
    @{ Cell = 'D34'; Value = '1.944' }
    @{ Cell = 'E34'; Value = '  +0.25%  ' }
    @{ Cell = 'D35'; Value = '0.7663' }
    @{ Cell = 'E35'; Value = '  +2.56%  ' }
    @{ Cell = 'E36'; Value = '  +0.43%  ' }
    @{ Cell = 'D37'; Value = '2.689' }
    @{ Cell = 'E37'; Value = '  -0.53%  ' }
    @{ Cell = 'D38'; Value = '0.01881' }
    @{ Cell = 'E38'; Value = '  +0.53%  ' }
    @{ Cell = 'D39'; Value = '1.255.26' }
    @{ Cell = 'E39'; Value = '  +2.30%  ' }
    @{ Cell = 'D40'; Value = '2.753' }
    @{ Cell = 'D41'; Value = '6.500' }
    @{ Cell = 'E41'; Value = '  -0.62%  ' }
    @{ Cell = 'D42'; Value = '0.9196' }
    @{ Cell = 'E42'; Value = '  +3.62%  ' }
    @{ Cell = 'D43'; Value = '112.82' }
    @{ Cell = 'E43'; Value = '  +2.08%  ' }
    @{ Cell = 'D44'; Value = '74.21' }
    @{ Cell = 'E44'; Value = '  +2.09%  ' }
    @{ Cell = 'E45'; Value = '  +8.54%  ' }
    @{ Cell = 'E46'; Value = '  +0.04%  ' }
    @{ Cell = 'D47'; Value = '2.048.90' }
    @{ Cell = 'E47'; Value = '  +1.42%  ' }
    @{ Cell = 'D48'; Value = '0.5224' }
    @{ Cell = 'E48'; Value = '  +0.57%  ' }
    @{ Cell = 'E49'; Value = '  -0.09%  ' }
    @{ Cell = 'D50'; Value = '9.465' }
    @{ Cell = 'E50'; Value = '  +0.79%  ' }
    @{ Cell = 'D51'; Value = '0.4363' }
    @{ Cell = 'E51'; Value = '  +1.03%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text storage so values like '25.00' or '6.500' keep their
    # exact textual representation instead of being coerced to numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.NumberFormat = "General"
}

